$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing timesheet rows 9, 10, 11 (copy date/time formatting from row 8) ---

# Row 9: 2014-11-16 (serial 41959)
$ws.Range("A8:C8").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A9").Value = 41959
$ws.Range("B9").Value = 0.916666666666667
$ws.Range("C9").Value = 0.999988425925926

# Row 10: 2014-11-17 (serial 41960)
$ws.Range("A8:C8").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A10").Value = 41960
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0.0833333333333333

# Row 11: 2014-11-18 (serial 41961)
$ws.Range("A8:C8").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A11").Value = 41961
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0.145833333333333

$excel.CutCopyMode = 0

# --- Add a grand total row summing the Total Time column ---
$ws.Range("D26").Formula = "=SUM(D3:D25)"

# --- Update the active selection to D10 ---
$ws.Range("D10").Select()
